# Bill of Materials AK70 - Hardware update
# 1. Rename "Chain (2ft)" -> "05B Chain (2ft)" (row 6)
# 2. Add new "Bowden Tube" line item (row 20) with link, count, price

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update chain description to specify 05B chain ---
$ws.Range("A6").Value = "05B Chain (2ft)"

# --- 2. Add new row 20: Bowden Tube ---
$bowdenUrl = "https://www.amazon.com/Jagwire-Sport-Housing-Slick-Lube-Titanium/dp/B085NBZMJS/ref=sr_1_1?crid=1A5WX5ADQYM0Y&dib=eyJ2IjoiMSJ9.CWV7EelBoN67bHKqG_VMGaAHwKq3lWAqkaCYopLdT43GyH4CDNeyWoQ_bFV_YrbZmhmwmsofP69GRzCCWYW_ULIkinZgZrdky8EGo_FPRa2GDLIPcrjwSu8T1nDFsZ03wyuLLatilRsdmpFkqvecV8S7AOhy1XjGVW6Ztcl1kgVL4_2zjOWbooP9z_kj4elJxdrMd7yL-uhr5ZcYm6F_Z725qnNy_c0-wNfaDt_xxhvLSabJzp2Ta9HzGpozVxtWGMo89NRT0qD_9iY1NxHHBedJp86w9TXZ0OfMm9s45WQ.7sQYhAOJzDOGkHvrSoGgv6F0CIzum6_7bjb1Zjv8zzc&dib_tag=se&keywords=Jagwire+Brake+Housing+CGX-SL+Slick-Lube+5+mm+%2810+m%29&qid=1726249759&s=sporting-goods&sprefix=jagwire+brake+housing+cgx-sl+slick-lube+5+mm+10+m+%2Csporting%2C99&sr=1-1"

$ws.Range("A20").Value = "Bowden Tube"
$ws.Range("B20").Value = $bowdenUrl
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 61.97

# Match the centered formatting used by the rest of the data rows (style index 5:
# default font/fill/border with horizontal-center alignment). Column B keeps the
# workbook's default (unstyled) formatting, same as every other Link cell.
$ws.Range("A20").HorizontalAlignment = -4108
$ws.Range("C20").HorizontalAlignment = -4108
$ws.Range("D20").HorizontalAlignment = -4108

# Reposition the window / selection like Excel would leave it after typing into A20:D20
$win = $wb.Windows.Item(1)
$win.Left = 28680
$win.Top = -120

$ws.Range("A20:D20").Select()
